$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "uODC"
$ws.Range("B27").Value = 6
$ws.Range("C27").Value = 12

$ws.Range("A28").Value = "BlackBox"
$ws.Range("B28").Value = 3
$ws.Range("C28").Value = 12

$ws.Range("A29").Value = "BlackBox - 10269269"
$ws.Range("B29").Value = 3
$ws.Range("C29").Value = 12

$ws.Range("A30").Value = "VertHor - 10176650"
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = 24
